# Air_Amadeus.xlsx - "code change for analysing failure test report"
#
# The Multicity sheet's test row gets a 4th leg appended to both the
# AirPortPairs and TravelDates columns, and a second scenario row is added
# that reuses the same (updated) test data but with IncludeNearByAirPorts
# turned on. The Multicity tab also becomes the active tab/sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Air_Amadeus_OneWay")
$ws3 = $wb.Worksheets.Item("Air_Amadeus_Multicity")

# 1. Extend the MultiCity itinerary (row 2) with a 4th leg: MIA-DEL / 80.
$ws3.Range("D2").Value = "LAS-LAX|LAX-DFW|DFW-MIA|MIA-DEL"
$ws3.Range("E2").Value = "53|64|71|80"

# 2. Add a new scenario row (row 3) under the table, copying the last data
#    row's table formatting (same style pattern Excel applies to the new
#    bottom row of an existing table) and then filling in the same test
#    values as row 2, except IncludeNearByAirPorts is turned on.
$ws1.Range("A3:Q3").Copy($ws3.Range("A3:Q3"))

$ws3.Range("A3").Value = "AmadeusWS air MULTICITY booking for DOMESTIC location for 1 Adult with Login."
$ws3.Range("C3").Value = "MultiCity"
$ws3.Range("D3").Value = "LAS-LAX|LAX-DFW|DFW-MIA|MIA-DEL"
$ws3.Range("E3").Value = "53|64|71|80"
$ws3.Range("F3").Value = 1
$ws3.Range("G3").Value = 0
$ws3.Range("H3").Value = 1
$ws3.Range("M3").Value = "Creditcard"
$ws3.Range("O3").Value = "Registered"

# 2b. Grow the worksheet Table (Table3) so it covers the new row too.
$tbl3 = $ws3.ListObjects.Item(1)
$tbl3.Resize($ws3.Range("A1:Q3"))

# 3. Make the MultiCity sheet the active tab/sheet, with E11 selected there.
$ws3.Activate()
$ws3.Range("E11").Select()
